# Update "想去人数" (want-to-go count) values in the 展览 sheet and the
# 全部类型 (aggregate) sheet to reflect newer scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): rows 3-10 in column F ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 258
$wsExhibit.Range("F4").Value = 135
$wsExhibit.Range("F5").Value = 1709
$wsExhibit.Range("F6").Value = 1495
$wsExhibit.Range("F7").Value = 274
$wsExhibit.Range("F8").Value = 61
$wsExhibit.Range("F9").Value = 455
$wsExhibit.Range("F10").Value = 119

# --- Sheet "全部类型" (sheet4): rows 3-7, 9-11 in column F (row 8 is an
#     unrelated concert entry and stays untouched) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 258
$wsAll.Range("F4").Value = 135
$wsAll.Range("F5").Value = 1709
$wsAll.Range("F6").Value = 1495
$wsAll.Range("F7").Value = 274
$wsAll.Range("F9").Value = 61
$wsAll.Range("F10").Value = 455
$wsAll.Range("F11").Value = 119
